$d = $word.ActiveDocument

# --- Grammar fix: "least lines of code" -> "fewest lines of code" ---
# Locate the word "least" that needs to become "fewest".
$rng = $d.Content
$rng.Find.Execute("least", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$leastStart = $rng.Start
$leastEnd = $rng.End

# Drop a throwaway bookmark exactly at the start of "least" first. This forces
# the run containing "least" to be split off from its (identically formatted)
# preceding sibling run, so replacing the word's text doesn't fold the two
# runs together.
$splitPoint = $d.Range($leastStart, $leastStart)
$d.Bookmarks.Add("ZZTempSplit", $splitPoint)

# Replace "least" with "fewest" in its now-isolated run.
$replacement = "fewest"
$editRng = $d.Range($leastStart, $leastEnd)
$editRng.Text = $replacement

$d.Bookmarks("ZZTempSplit").Delete()

# --- Move the _GoBack bookmark so it sits right after "fewest" ---
# (it previously sat immediately after "Docker Swarm" in the next paragraph).
$d.Bookmarks("_GoBack").Delete()

$fewestEnd = $leastStart + $replacement.Length
$goBackRange = $d.Range($fewestEnd, $fewestEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange)
